$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting the existing A:D columns to B:E.
$ws.Range("A1:A15").EntireColumn.Insert()

# Update header row text (now shifted into D1:E1) to use "n = " spacing.
$ws.Range("D1").Value = "Treatment at T1 (n = 5080)"
$ws.Range("E1").Value = "Control at T1 (n = 745)"

# Update the "P=" labels (now in column B) to use "P = " spacing.
$ws.Range("B3").Value = "Gender (P = 0.006)"
$ws.Range("B10").Value = "Interested in News (P = 0.000)"

# Give the now-empty column A (rows 2-15) the same header style as row 1
# (bold, bordered, centered) by copying the format from B1.
$ws.Range("B1").Copy()
$ws.Range("A2:A15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
